$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3239.8
$ws.Range("I62").Value = 3239.8
$ws.Range("K62").Value = 3239.8
$ws.Range("M62").Value = -2615.8
$ws.Range("H65").Value = 3239.8
$ws.Range("I65").Value = 3239.8
$ws.Range("K65").Value = 16199
$ws.Range("M65").Value = -13079
$ws.Range("I103").Value = 452.33334
$ws.Range("J103").Value = 62501984
$ws.Range("K103").Value = 1357.00002
$ws.Range("L103").Value = 187505952
$ws.Range("M103").Value = -771.0000199999999
$ws.Range("N103").Value = -187507124
$ws.Range("H137").Value = 2675.8064
$ws.Range("J137").Value = 4142.4546
$ws.Range("L137").Value = 12427.3638
$ws.Range("N137").Value = -17527.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7709.757
$ws.Range("I32").Value = 7472.9043
$ws.Range("K32").Value = 7472.9043
$ws.Range("M32").Value = -7185.9043
$ws.Range("H97").Value = 1437.0968
$ws.Range("I97").Value = 818.2105
$ws.Range("J97").Value = 2417
$ws.Range("K97").Value = 818.2105
$ws.Range("L97").Value = 2417
$ws.Range("M97").Value = -322.2105
$ws.Range("N97").Value = -3409

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5302.343
$ws.Range("I20").Value = 5447.185
$ws.Range("K20").Value = 5447.185
$ws.Range("M20").Value = -5200.185
$ws.Range("H86").Value = 787845.0600000001
$ws.Range("I86").Value = 1001475.8
$ws.Range("J86").Value = 4532.3335
$ws.Range("K86").Value = 1001475.8
$ws.Range("L86").Value = 4532.3335
$ws.Range("M86").Value = -1000352.8
$ws.Range("N86").Value = -6778.3335
$ws.Range("H89").Value = 787845.0600000001
$ws.Range("I89").Value = 1001475.8
$ws.Range("J89").Value = 4532.3335
$ws.Range("K89").Value = 5007379
$ws.Range("L89").Value = 22661.6675
$ws.Range("M89").Value = -5001763
$ws.Range("N89").Value = -33893.6675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41670404
$ws.Range("I31").Value = 55558020
$ws.Range("K31").Value = 55558020
$ws.Range("M31").Value = -55557725
$ws.Range("H34").Value = 41670404
$ws.Range("I34").Value = 55558020
$ws.Range("K34").Value = 55558020
$ws.Range("M34").Value = -55557818
$ws.Range("H141").Value = 421817.72
$ws.Range("I141").Value = 150000
$ws.Range("K141").Value = 150000
$ws.Range("M141").Value = -144820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 5992.091
$ws.Range("I60").Value = 299.57144
$ws.Range("J60").Value = 15954
$ws.Range("K60").Value = 898.71432
$ws.Range("L60").Value = 47862
$ws.Range("M60").Value = -647.71432
$ws.Range("N60").Value = -48364
$ws.Range("H107").Value = 3640707.8
$ws.Range("I107").Value = 4712
$ws.Range("J107").Value = 4333278.5
$ws.Range("K107").Value = 14136
$ws.Range("L107").Value = 12999835.5
$ws.Range("M107").Value = -12216
$ws.Range("N107").Value = -13003675.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 15501.5
$ws.Range("J10").Value = 15501.5
$ws.Range("L10").Value = 15501.5
$ws.Range("N10").Value = -15839.5
$ws.Range("H35").Value = 43722
$ws.Range("I35").Value = 43000
$ws.Range("K35").Value = 43000
$ws.Range("M35").Value = -42702
$ws.Range("H52").Value = 23665.666
$ws.Range("I52").Value = 23665.666
$ws.Range("K52").Value = 23665.666
$ws.Range("M52").Value = -23406.666
$ws.Range("H136").Value = 6471.7617
$ws.Range("J136").Value = 6471.7617
$ws.Range("L136").Value = 19415.2851
$ws.Range("N136").Value = -24515.2851

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8955.137000000001
$ws.Range("I7").Value = 9831.182000000001
$ws.Range("J7").Value = 8079.091
$ws.Range("K7").Value = 9831.182000000001
$ws.Range("L7").Value = 8079.091
$ws.Range("M7").Value = -9719.182000000001
$ws.Range("N7").Value = -8303.091
$ws.Range("H40").Value = 5756.7085
$ws.Range("I40").Value = 4509.278
$ws.Range("K40").Value = 4509.278
$ws.Range("M40").Value = -4373.278
$ws.Range("H76").Value = 47999
$ws.Range("J76").Value = 47999
$ws.Range("L76").Value = 47999
$ws.Range("N76").Value = -48675
$ws.Range("H79").Value = 47999
$ws.Range("J79").Value = 47999
$ws.Range("L79").Value = 47999
$ws.Range("N79").Value = -50339
$ws.Range("H122").Value = 3753.8542
$ws.Range("I122").Value = 3482.3262
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 10446.9786
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -7996.9786
$ws.Range("N122").Value = -34897
$ws.Range("H126").Value = 8955.137000000001
$ws.Range("I126").Value = 9831.182000000001
$ws.Range("J126").Value = 8079.091
$ws.Range("K126").Value = 29493.546
$ws.Range("L126").Value = 24237.273
$ws.Range("M126").Value = -27023.546
$ws.Range("N126").Value = -29177.273
$ws.Range("H128").Value = 69995
$ws.Range("J128").Value = 69995
$ws.Range("L128").Value = 69995
$ws.Range("N128").Value = -79955
$ws.Range("H137").Value = 116985.14
$ws.Range("J137").Value = 116985.14
$ws.Range("L137").Value = 116985.14
$ws.Range("N137").Value = -127185.14

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 34498
$ws.Range("J63").Value = 34498
$ws.Range("L63").Value = 34498
$ws.Range("N63").Value = -35746
$ws.Range("H66").Value = 34498
$ws.Range("J66").Value = 34498
$ws.Range("L66").Value = 103494
$ws.Range("N66").Value = -109734
$ws.Range("H113").Value = 630.28125
$ws.Range("I113").Value = 526.0417
$ws.Range("J113").Value = 943
$ws.Range("K113").Value = 1578.1251
$ws.Range("L113").Value = 2829
$ws.Range("M113").Value = 591.8749
$ws.Range("N113").Value = -7169
$ws.Range("H122").Value = 2684.7222
$ws.Range("I122").Value = 2288.3333
$ws.Range("K122").Value = 6864.999899999999
$ws.Range("M122").Value = -4414.999899999999
$ws.Range("H126").Value = 4396.524
$ws.Range("I126").Value = 4019.2942
$ws.Range("K126").Value = 12057.8826
$ws.Range("M126").Value = -9587.882599999999
$ws.Range("H132").Value = 252474.33
$ws.Range("I132").Value = 2315.5806
$ws.Range("J132").Value = 1114132.2
$ws.Range("K132").Value = 6946.7418
$ws.Range("L132").Value = 3342396.6
$ws.Range("M132").Value = -4416.7418
$ws.Range("N132").Value = -3347456.6
